$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Extend the header row (row 1) with two new columns: P1=14, Q1=15 ---
# Copy formatting from the last existing header cell (O1) so the new
# header cells pick up the same bold/centered/bordered style without
# minting new style records.
$ws.Range("O1").Copy() | Out-Null
$ws.Range("P1:Q1").PasteSpecial(-4122) | Out-Null
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# --- Update the data rows (2-25) ---
# Columns I and K swap their values (1<->2), and columns M and O swap
# their values (1<->2) as well; then two new trailing columns P and Q
# are appended, each holding the value 2.
for ($r = 2; $r -le 25; $r++) {
    $ws.Range("I$r").Value = 2
    $ws.Range("K$r").Value = 1
    $ws.Range("M$r").Value = 2
    $ws.Range("O$r").Value = 1
    $ws.Range("P$r").Value = 2
    $ws.Range("Q$r").Value = 2
}

$excel.CutCopyMode = 0

Write-Output "applied parallel.xlsx contingency edit (A1:O25 -> A1:Q25)"
